# Fixing name of Sectors to be alligned with Baseline
# Rename the four commodity/sector header labels (row 3, columns D:G) on every
# year sheet (2000 .. 2100) from their short codes to their full names.

$wb = $excel.ActiveWorkbook

$sectorAddrs = @("D3", "E3", "F3", "G3")
$sectorNames = @("Neodymium", "Dysprosium", "Copper ores and concentrates", "Raw silicon")

foreach ($ws in $wb.Worksheets) {
    for ($i = 0; $i -lt $sectorAddrs.Length; $i++) {
        $ws.Range($sectorAddrs[$i]).Value = $sectorNames[$i]
    }
}

# Tiny last-digit recalculation drift on the "Wires" total (row 7, col G) that
# came along with the baseline re-run, for the affected year sheets.
$g7Years = @("2018", "2027", "2029", "2041", "2042", "2047", "2055", "2058", "2069", "2072", "2073", "2081", "2097", "2099", "2100")
$g7Values = @(-818895.9589051851, -68343063.31331737, -89179111.04986022, -334085716.730737, -363748238.1360357, -847938034.7792233, -2669462722.250011, -2637391018.485797, -333862256.0526412, -239084158.9771962, -250643173.2307304, -255985131.7648519, -290428733.8834999, -300597185.8741556, -296289626.1066293)

for ($i = 0; $i -lt $g7Years.Length; $i++) {
    $ws = $wb.Worksheets.Item($g7Years[$i])
    $ws.Range("G7").Value = $g7Values[$i]
}
